# Refresh the cryptocurrency price list (columns B:Coin, C:Link,
# D:Price, E:Volume(1h)) to the latest scraped snapshot.
#
# Column D holds "Price" as plain text even when it looks numeric
# (e.g. "574.03", or "66.241.84" which uses '.' as a thousands
# separator and is therefore never a valid number anyway). Writing a
# numeric-looking string straight into Range.Value lets Excel's usual
# text-to-number inference kick in and silently convert/round it, so
# for those cells we briefly force Text number-formatting, assign the
# value, then clear the formatting again so the cell keeps its
# original (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "66.338.71"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "3.242.26"
$ws.Range("E3").Value = "  +4.85%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "576.45"
$ws.Range("E5").Value = "  +1.66%  "
Set-TextValue "D6" "155.43"
$ws.Range("E6").Value = "  +8.17%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.232.05"
$ws.Range("E8").Value = "  +4.74%  "
Set-TextValue "D9" "0.517"
$ws.Range("E9").Value = "  +3.87%  "
$ws.Range("E10").Value = "  +10.52%  "
Set-TextValue "D11" "0.168"
$ws.Range("E11").Value = "  +5.62%  "
Set-TextValue "D12" "0.487"
$ws.Range("E12").Value = "  +4.01%  "
Set-TextValue "D13" "38.41"
$ws.Range("E13").Value = "  +6.75%  "
Set-TextValue "D14" "0.0000237"
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("D15").Value = "3.762.25"
$ws.Range("E15").Value = "  +5.13%  "
$ws.Range("D16").Value = "66.326.80"
$ws.Range("E16").Value = "  +2.49%  "
Set-TextValue "D17" "548.46"
$ws.Range("E17").Value = "  +10.19%  "
$ws.Range("D18").Value = "3.247.36"
$ws.Range("E18").Value = "  +5.24%  "
$ws.Range("E19").Value = "  +3.13%  "
Set-TextValue "D20" "7.11"
$ws.Range("E20").Value = "  +6.11%  "
Set-TextValue "D21" "14.55"
$ws.Range("E21").Value = "  +5.13%  "
Set-TextValue "D22" "0.743"
$ws.Range("E22").Value = "  +7.07%  "
Set-TextValue "D23" "7.80"
$ws.Range("E23").Value = "  +7.85%  "
Set-TextValue "D24" "13.57"
$ws.Range("E24").Value = "  +6.31%  "
Set-TextValue "D25" "82.26"
$ws.Range("E25").Value = "  +4.24%  "
$ws.Range("E26").Value = "  +0.03%  "
Set-TextValue "D27" "9.38"
$ws.Range("E27").Value = "  +15.40%  "
Set-TextValue "D28" "2.91"
$ws.Range("E28").Value = "  +4.62%  "
Set-TextValue "D29" "2.28"
$ws.Range("E29").Value = "  +8.73%  "
Set-TextValue "D30" "27.94"
$ws.Range("E30").Value = "  +4.95%  "
Set-TextValue "D31" "2.78"
$ws.Range("E31").Value = "  +3.29%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +4.79%  "
Set-TextValue "D34" "573.46"
$ws.Range("E34").Value = "  +10.37%  "
Set-TextValue "D35" "5.85"
$ws.Range("E35").Value = "  +4.74%  "
Set-TextValue "D36" "6.48"
$ws.Range("E36").Value = "  +7.78%  "
Set-TextValue "D37" "0.0471"
$ws.Range("E37").Value = "  +14.74%  "
Set-TextValue "D38" "54.93"
$ws.Range("E38").Value = "  +3.15%  "
Set-TextValue "D39" "0.0877"
$ws.Range("E39").Value = "  +9.45%  "
$ws.Range("E40").Value = "  +13.86%  "
$ws.Range("E41").Value = "  +4.38%  "
$ws.Range("D42").Value = "3.144.29"
$ws.Range("E42").Value = "  +6.21%  "
Set-TextValue "D43" "8.68"
$ws.Range("E43").Value = "  +3.11%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D44" "0.276"
$ws.Range("E44").Value = "  +11.28%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D45" "2.36"
$ws.Range("E45").Value = "  +9.73%  "
Set-TextValue "D46" "27.26"
$ws.Range("E46").Value = "  +7.31%  "
$ws.Range("D47").Value = "0.0₃0572"
$ws.Range("E47").Value = "  +4.35%  "
Set-TextValue "D49" "0.114"
$ws.Range("E49").Value = "  +4.63%  "
Set-TextValue "D50" "2.27"
$ws.Range("E50").Value = "  +8.42%  "
Set-TextValue "D51" "122.68"
$ws.Range("E51").Value = "  +1.13%  "
